$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (B and C) before the existing STATUS column,
# pushing the original "STATUS" header from B1 to D1 while keeping its
# formatting (border/bold/alignment style).
$ws.Columns("B:C").Insert()

# New header cells - they inherit the header style carried over by the
# column insert, same as the original "E-mail"/"STATUS" headers.
$ws.Range("B1").Value = "Nome"
$ws.Range("C1").Value = "Produto"

# New data rows with the e-mail send results.
$ws.Range("A2").Value = "oversouls11@gmail.com"
$ws.Range("D2").Value = "SUCESSO"

$ws.Range("A3").Value = "davinascimento860.ld@gmail.com"

# Restore selection to A1 (matches the saved workbook view state).
$ws.Range("A1").Select()
